# Regenerate the attribute-name header row (C1:N1) on Sheet1.
#
# The source data table was regenerated, which reshuffled the order the
# attribute-code columns (ARC, BC, CM, FM, FOC, IV, LR, OV, POC, SV, SOB, VC)
# are laid out in. The new left-to-right header order is:
#   SV, ARC, LR, FM, FOC, POC, OV, BC, IV, VC, CM, SOB

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C1").Value = "SV"
$ws.Range("D1").Value = "ARC"
$ws.Range("E1").Value = "LR"
$ws.Range("F1").Value = "FM"
$ws.Range("G1").Value = "FOC"
$ws.Range("H1").Value = "POC"
$ws.Range("I1").Value = "OV"
$ws.Range("J1").Value = "BC"
$ws.Range("K1").Value = "IV"
$ws.Range("L1").Value = "VC"
$ws.Range("M1").Value = "CM"
$ws.Range("N1").Value = "SOB"

# Re-apply "best fit" column widths for the re-labelled columns (the new
# header text has different character widths than the old one).
$ws.Columns.Item(3).ColumnWidth  = 3.0              # C  SV
$ws.Columns.Item(4).ColumnWidth  = 4.428571428571429 # D  ARC
$ws.Columns.Item(5).ColumnWidth  = 2.857142857142857 # E  LR
$ws.Columns.Item(6).ColumnWidth  = 3.571428571428572 # F  FM
$ws.Columns.Item(7).ColumnWidth  = 4.571428571428571 # G  FOC
$ws.Columns.Item(8).ColumnWidth  = 4.714285714285714 # H  POC
$ws.Columns.Item(9).ColumnWidth  = 3.571428571428572 # I  OV
$ws.Columns.Item(10).ColumnWidth = 3.142857142857143 # J  BC
$ws.Columns.Item(11).ColumnWidth = 2.571428571428572 # K  IV
$ws.Columns.Item(12).ColumnWidth = 3.285714285714286 # L  VC
$ws.Columns.Item(13).ColumnWidth = 3.857142857142857 # M  CM
$ws.Columns.Item(14).ColumnWidth = 4.571428571428571 # N  SOB
